# "moving ref params earlier"
#
# The last four slides of the deck (sldId 448, 450, 451, 452 -- the
# "Implementing Reference parameters" / "A change to extend-env" slides
# and friends) are being removed from this location. Deleting them also
# drops the notes slide that was only referenced by one of them.
#
# Delete from the end of the deck backwards so slide indices stay valid
# as each Delete() call shifts everything after it.

$p = $ppt.ActivePresentation

$p.Slides.Item(33).Delete()
$p.Slides.Item(32).Delete()
$p.Slides.Item(31).Delete()
$p.Slides.Item(30).Delete()
